$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    # Force the cell to remain text even when $val looks numeric
    # (the source file stores these as plain text/inlineStr cells).
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "26.065.90"
$ws.Range("E2").Value = "  -0.10%  "

Set-TextValue "D3" "1.638.69"
$ws.Range("E3").Value = "  -1.62%  "

$ws.Range("E4").Value = "  -0.19%  "

Set-TextValue "D5" "213.77"
$ws.Range("E5").Value = "  +2.00%  "

Set-TextValue "D6" "0.5239"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("E8").Value = "  -0.99%  "

Set-TextValue "D9" "0.06303"
$ws.Range("E9").Value = "  +0.28%  "

Set-TextValue "D10" "20.69"
$ws.Range("E10").Value = "  -2.05%  "

Set-TextValue "D11" "0.07660"
$ws.Range("E11").Value = "  +1.63%  "

Set-TextValue "D12" "1.633.10"
$ws.Range("E12").Value = "  -1.97%  "

$ws.Range("E13").Value = "  -0.43%  "

Set-TextValue "D14" "1.861.11"
$ws.Range("E14").Value = "  -1.74%  "

Set-TextValue "D15" "0.5526"
$ws.Range("E15").Value = "  +0.06%  "

Set-TextValue "D16" "0.0₅8285"
$ws.Range("E16").Value = "  +4.74%  "

$ws.Range("E17").Value = "  -2.31%  "

Set-TextValue "D18" "26.045.56"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("E19").Value = "  -0.11%  "

Set-TextValue "D20" "4.693"
$ws.Range("E20").Value = "  -0.58%  "

Set-TextValue "D21" "188.14"
$ws.Range("E21").Value = "  +0.96%  "

Set-TextValue "D22" "10.19"
$ws.Range("E22").Value = "  -1.05%  "

Set-TextValue "D23" "6.162"
$ws.Range("E23").Value = "  +0.15%  "

Set-TextValue "D25" "145.54"
$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("E26").Value = "  -2.28%  "

Set-TextValue "D27" "7.421"
$ws.Range("E27").Value = "  -0.66%  "

Set-TextValue "D28" "15.80"
$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("E29").Value = "  +2.98%  "

Set-TextValue "D30" "0.05967"
$ws.Range("E30").Value = "  -5.37%  "

Set-TextValue "D31" "1.255"
$ws.Range("E31").Value = "  -1.38%  "

Set-TextValue "D32" "3.440"
$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("E34").Value = "  +0.71%  "

Set-TextValue "D35" "0.9848"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D36" "2.761"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D37" "2.393"
$ws.Range("E37").Value = "  -0.57%  "

Set-TextValue "D38" "0.5681"
$ws.Range("E38").Value = "  -5.60%  "

$ws.Range("E39").Value = "  +0.24%  "

Set-TextValue "D40" "0.8492"
$ws.Range("E40").Value = "  -2.33%  "

Set-TextValue "D41" "5.741"
$ws.Range("E41").Value = "  -5.52%  "

Set-TextValue "D42" "1.000"
$ws.Range("E42").Value = "  -0.26%  "

Set-TextValue "D43" "1.034.28"
$ws.Range("E43").Value = "  -6.51%  "

Set-TextValue "D44" "100.21"
$ws.Range("E44").Value = "  +0.49%  "

Set-TextValue "D45" "1.786.76"
$ws.Range("E45").Value = "  -1.69%  "

Set-TextValue "D46" "0.0₈106"
$ws.Range("E46").Value = "  -2.70%  "

Set-TextValue "D47" "55.78"
$ws.Range("E47").Value = "  +0.84%  "

Set-TextValue "D48" "0.9952"
$ws.Range("E48").Value = "  -0.86%  "

Set-TextValue "D49" "8.052"
$ws.Range("E49").Value = "  +0.62%  "

Set-TextValue "D50" "0.05156"
$ws.Range("E50").Value = "  -1.53%  "

Set-TextValue "D51" "0.4214"
$ws.Range("E51").Value = "  -0.67%  "
